$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 88

$ws.Cells.Item($newRow, 1).Value = "CompaNanny"
$ws.Cells.Item($newRow, 2).Value = "CompaNanny Amstelveen KDV"
$ws.Cells.Item($newRow, 3).Value = "KDV"
# Leading apostrophe keeps this a literal text value ("2024-02-07")
# instead of Excel auto-converting it into a date serial number,
# matching the other rows in column D.
$ws.Cells.Item($newRow, 4).Value = "'2024-02-07"
$ws.Cells.Item($newRow, 5).Value = 0
$ws.Cells.Item($newRow, 6).Value = 0
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
